# Generate Report for Handback
# - Overview: row for the 31028f7b... (de-de) handoff moves from
#   "Ready for handoff" to "Handback transform failed" (Status + Source Status cols).
# - zh-cn sheet: record the handback/handoff filename mismatch error in the
#   "Error Detail" column (K) for the 31028f7b... row (row 3).
# - de-de sheet: same error, recorded on its own sheet's row 3.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

$newStatus = "Handback transform failed"

# Every cell that previously displayed "Ready for handoff" (the Status of
# the 31028f7b... handoff, row 3 on each sheet) now reads "Handback
# transform failed" - update them all so they continue to share the text.
$overview.Range("B3").Value = $newStatus
$overview.Range("C3").Value = $newStatus
$overview.Range("D3").Value = $newStatus
$zhcn.Range("C3").Value = $newStatus
$dede.Range("C3").Value = $newStatus

$zhcnError = "Handback file name: f3c2ifpa.2y2 is different with handoff file name: 31028f7b-8c46-4933-b307-2043f0d0d677.547f760b68d8ef99276b6c74ff34bd946778d27b.zh-cn."
$dedeError = "Handback file name: f3c2ifpa.2y2 is different with handoff file name: 31028f7b-8c46-4933-b307-2043f0d0d677.547f760b68d8ef99276b6c74ff34bd946778d27b.de-de."

# Column K = "Error Detail" on both language sheets; row 3 = the
# 31028f7b-8c46-4933-b307-2043f0d0d677 handoff.
$zhcn.Range("K3").Value = $zhcnError
$dede.Range("K3").Value = $dedeError
